# Update "想去人数" (want-to-go count) figures by +1 for a handful of
# events across the "展览" (Exhibitions), "本地生活" (Local Life) and
# "全部类型" (All Types) sheets, matching the refreshed scrape snapshot.

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions) sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 902
$ws1.Range("F6").Value = 331
$ws1.Range("F7").Value = 673
$ws1.Range("F30").Value = 49
$ws1.Range("F31").Value = 100
$ws1.Range("F38").Value = 15
$ws1.Range("F46").Value = 79

# 本地生活 (Local Life) sheet
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 761

# 全部类型 (All Types) sheet — aggregated view of the above rows
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 761
$ws4.Range("F7").Value = 902
$ws4.Range("F10").Value = 331
$ws4.Range("F11").Value = 673
$ws4.Range("F32").Value = 49
$ws4.Range("F33").Value = 100
